# Generate Report for Handoff
# Updates the Priority and handoff timestamps for the rows that were
# previously marked "low" priority (the 02e525cb... file group), to reflect
# the newly generated handoff report ("ht" priority, refreshed timestamps).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 4-7
for ($r = 4; $r -le 7; $r++) {
    $overview.Range("G$r").Value = "2016-08-20 14:36:30"
}

# zh-cn sheet: Priority column (E) flips from "low" to "ht", and the
# "Latest Handoff Datetime" column (H) is refreshed, rows 4-7
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-20 14:36:26"
}

# de-de sheet: Priority column (E) flips from "low" to "ht", and the
# "Latest Handoff Datetime" column (H) is refreshed, rows 4-7
for ($r = 4; $r -le 7; $r++) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-20 14:36:30"
}
